$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 6.206015333333333
$ws.Range("H2").Value = 18.618046
$ws.Range("I2").Value = 0.0150172404156507
$ws.Range("J2").Value = 0.0150172404156507
$ws.Range("M2").Value = 0.6415476666666667
$ws.Range("N2").Value = 1.924643
$ws.Range("O2").Value = 0.1426849042655057
$ws.Range("P2").Value = 0.1426849042655057
$ws.Range("Q2").Value = 3.981454656397556
$ws.Range("R2").Value = 35.833091907578
$ws.Range("S2").Value = 0.002142733511039204
$ws.Range("T2").Value = 0.002142733511039203

# Row 3
$ws.Range("G3").Value = 6.206015333333333
$ws.Range("H3").Value = 18.618046
$ws.Range("I3").Value = 0.0150172404156507
$ws.Range("J3").Value = 0.0150172404156507
$ws.Range("O3").Value = 0.5986102210699216
$ws.Range("P3").Value = 0.5986102210699217
$ws.Range("Q3").Value = 16.703515093728
$ws.Range("R3").Value = 150.331635843552
$ws.Range("S3").Value = 0.008989473605072826
$ws.Range("T3").Value = 0.008989473605072826

# Row 4
$ws.Range("G4").Value = 6.206015333333333
$ws.Range("H4").Value = 18.618046
$ws.Range("I4").Value = 0.0150172404156507
$ws.Range("J4").Value = 0.0150172404156507
$ws.Range("M4").Value = 1.163203
$ws.Range("N4").Value = 3.489609
$ws.Range("O4").Value = 0.2587048746645726
$ws.Range("P4").Value = 0.2587048746645726
$ws.Range("Q4").Value = 7.218855653779333
$ws.Range("R4").Value = 64.96970088401399
$ws.Range("S4").Value = 0.003885033299538669
$ws.Range("T4").Value = 0.003885033299538669

# Row 5
$ws.Range("I5").Value = 0.9317452840597572
$ws.Range("J5").Value = 0.9317452840597571
$ws.Range("M5").Value = 0.6415476666666667
$ws.Range("N5").Value = 1.924643
$ws.Range("O5").Value = 0.1426849042655057
$ws.Range("P5").Value = 0.1426849042655057
$ws.Range("Q5").Value = 247.0295138865859
$ws.Range("R5").Value = 2223.265624979274
$ws.Range("S5").Value = 0.1329459866559029
$ws.Range("T5").Value = 0.1329459866559029

# Row 6
$ws.Range("I6").Value = 0.9317452840597572
$ws.Range("J6").Value = 0.9317452840597571
$ws.Range("O6").Value = 0.5986102210699216
$ws.Range("P6").Value = 0.5986102210699217
$ws.Range("S6").Value = 0.5577522504718682
$ws.Range("T6").Value = 0.5577522504718682

# Row 7
$ws.Range("I7").Value = 0.9317452840597572
$ws.Range("J7").Value = 0.9317452840597571
$ws.Range("M7").Value = 1.163203
$ws.Range("N7").Value = 3.489609
$ws.Range("O7").Value = 0.2587048746645726
$ws.Range("P7").Value = 0.2587048746645726
$ws.Range("Q7").Value = 447.8941886491444
$ws.Range("R7").Value = 4031.047697842299
$ws.Range("S7").Value = 0.2410470469319861
$ws.Range("T7").Value = 0.2410470469319861

# Row 8
$ws.Range("G8").Value = 22.00088566666667
$ws.Range("H8").Value = 66.002657
$ws.Range("I8").Value = 0.05323747552459213
$ws.Range("J8").Value = 0.05323747552459213
$ws.Range("M8").Value = 0.6415476666666667
$ws.Range("N8").Value = 1.924643
$ws.Range("O8").Value = 0.1426849042655057
$ws.Range("P8").Value = 0.1426849042655057
$ws.Range("Q8").Value = 14.11461686405011
$ws.Range("R8").Value = 127.031551776451
$ws.Range("S8").Value = 0.007596184098563633
$ws.Range("T8").Value = 0.007596184098563632

# Row 9
$ws.Range("G9").Value = 22.00088566666667
$ws.Range("H9").Value = 66.002657
$ws.Range("I9").Value = 0.05323747552459213
$ws.Range("J9").Value = 0.05323747552459213
$ws.Range("O9").Value = 0.5986102210699216
$ws.Range("P9").Value = 0.5986102210699217
$ws.Range("Q9").Value = 59.215471775376
$ws.Range("R9").Value = 532.9392459783841
$ws.Range("S9").Value = 0.03186849699298064
$ws.Range("T9").Value = 0.03186849699298064

# Row 10
$ws.Range("G10").Value = 22.00088566666667
$ws.Range("H10").Value = 66.002657
$ws.Range("I10").Value = 0.05323747552459213
$ws.Range("J10").Value = 0.05323747552459213
$ws.Range("M10").Value = 1.163203
$ws.Range("N10").Value = 3.489609
$ws.Range("O10").Value = 0.2587048746645726
$ws.Range("P10").Value = 0.2587048746645726
$ws.Range("Q10").Value = 25.59149621012367
$ws.Range("R10").Value = 230.323465891113
$ws.Range("S10").Value = 0.01377279443304786
$ws.Range("T10").Value = 0.01377279443304786
